$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Tnfsf10", "Tnfrsf11b", "ECs", 3, 1, 15.77467133333333, 47.324014, 0.6679291979451994, 0.6679291979451996, 1, 0.3333333333333333, 0.07580833333333332, 0.227425, 0.01036056733223033, 0.01036056733223033, 1.195851542661111, 10.76266388395, 0.006920125428473839, 0.006920125428473841),
    @("ECs", "Tnfsf10", "Tnfrsf11b", "FAPs", 3, 1, 15.77467133333333, 47.324014, 0.6679291979451994, 0.6679291979451996, 3, 1, 6.93109, 20.79327, 0.9472576624920083, 0.9472576624920083, 109.3356667317533, 984.0210005857799, 0.6327010507557316, 0.6327010507557316),
    @("ECs", "Tnfsf10", "Tnfrsf11b", "sCs", 3, 1, 15.77467133333333, 47.324014, 0.6679291979451994, 0.6679291979451996, 3, 1, 0.3101076666666667, 0.930323, 0.04238177017576132, 0.04238177017576132, 4.891846519613555, 44.026618676522, 0.02830802176099403, 0.02830802176099404),
    @("FAPs", "Tnfsf10", "Tnfrsf11b", "ECs", 3, 1, 1.503082, 4.509246, 0.06364331360644089, 0.0636433136064409, 1, 0.3333333333333333, 0.07580833333333332, 0.227425, 0.01036056733223033, 0.01036056733223033, 0.1139461412833333, 1.02551527155, 0.0006593808358657815, 0.0006593808358657817),
    @("FAPs", "Tnfsf10", "Tnfrsf11b", "FAPs", 3, 1, 1.503082, 4.509246, 0.06364331360644089, 0.0636433136064409, 3, 1, 6.93109, 20.79327, 0.9472576624920083, 0.9472576624920083, 10.41799661938, 93.76196957442001, 0.06028661648008302, 0.06028661648008304),
    @("FAPs", "Tnfsf10", "Tnfrsf11b", "sCs", 3, 1, 1.503082, 4.509246, 0.06364331360644089, 0.0636433136064409, 3, 1, 0.3101076666666667, 0.930323, 0.04238177017576132, 0.04238177017576132, 0.4661172518286667, 4.195055266458001, 0.002697316290492081, 0.002697316290492082),
    @("M1", "Tnfsf10", "Tnfrsf11b", "ECs", 3, 1, 2.782582666666667, 8.347748000000001, 0.1178197738317093, 0.1178197738317093, 1, 0.3333333333333333, 0.07580833333333332, 0.227425, 0.01036056733223033, 0.01036056733223033, 0.2109429543222222, 1.8984865889, 0.001220679699851573, 0.001220679699851574),
    @("M1", "Tnfsf10", "Tnfrsf11b", "FAPs", 3, 1, 2.782582666666667, 8.347748000000001, 0.1178197738317093, 0.1178197738317093, 3, 1, 6.93109, 20.79327, 0.9472576624920083, 0.9472576624920083, 19.28633089510667, 173.57697805596, 0.111605683555162, 0.111605683555162),
    @("M1", "Tnfsf10", "Tnfrsf11b", "sCs", 3, 1, 2.782582666666667, 8.347748000000001, 0.1178197738317093, 0.1178197738317093, 3, 1, 0.3101076666666667, 0.930323, 0.04238177017576132, 0.04238177017576132, 0.8629002180671111, 7.766101962604001, 0.004993410576695681, 0.004993410576695682),
    @("M2", "Tnfsf10", "Tnfrsf11b", "ECs", 3, 1, 2.906420333333333, 8.719260999999999, 0.1230632931180533, 0.1230632931180533, 1, 0.3333333333333333, 0.07580833333333332, 0.227425, 0.01036056733223033, 0.01036056733223033, 0.2203308814361111, 1.982977932925, 0.001275005534475588, 0.001275005534475589),
    @("M2", "Tnfsf10", "Tnfrsf11b", "FAPs", 3, 1, 2.906420333333333, 8.719260999999999, 0.1230632931180533, 0.1230632931180533, 3, 1, 6.93109, 20.79327, 0.9472576624920083, 0.9472576624920083, 20.14466090816333, 181.30194817347, 0.116572647377576, 0.116572647377576),
    @("M2", "Tnfsf10", "Tnfrsf11b", "sCs", 3, 1, 2.906420333333333, 8.719260999999999, 0.1230632931180533, 0.1230632931180533, 3, 1, 0.3101076666666667, 0.930323, 0.04238177017576132, 0.04238177017576132, 0.9013032279225555, 8.111729051303, 0.005215640206001684, 0.005215640206001685),
    @("Neutro", "Tnfsf10", "Tnfrsf11b", "ECs", 3, 1, 0.3953666666666667, 1.1861, 0.01674056688603805, 0.01674056688603805, 1, 0.3333333333333333, 0.07580833333333332, 0.227425, 0.01036056733223033, 0.01036056733223033, 0.02997208805555555, 0.2697487925, 0.0001734417704025027, 0.0001734417704025027),
    @("Neutro", "Tnfsf10", "Tnfrsf11b", "FAPs", 3, 1, 0.3953666666666667, 1.1861, 0.01674056688603805, 0.01674056688603805, 3, 1, 6.93109, 20.79327, 0.9472576624920083, 0.9472576624920083, 2.740321949666667, 24.662897547, 0.01585763025725952, 0.01585763025725952),
    @("Neutro", "Tnfsf10", "Tnfrsf11b", "sCs", 3, 1, 0.3953666666666667, 1.1861, 0.01674056688603805, 0.01674056688603805, 3, 1, 0.3101076666666667, 0.930323, 0.04238177017576132, 0.04238177017576132, 0.1226062344777778, 1.1034561103, 0.0007094948583760252, 0.0007094948583760252),
    @("sCs", "Tnfsf10", "Tnfrsf11b", "ECs", 3, 1, 0.2551576666666667, 0.765473, 0.01080385461255898, 0.01080385461255898, 1, 0.3333333333333333, 0.07580833333333332, 0.227425, 0.01036056733223033, 0.01036056733223033, 0.01934307744722222, 0.174087697025, 0.0001119340631610445, 0.0001119340631610445),
    @("sCs", "Tnfsf10", "Tnfrsf11b", "FAPs", 3, 1, 0.2551576666666667, 0.765473, 0.01080385461255898, 0.01080385461255898, 3, 1, 6.93109, 20.79327, 0.9472576624920083, 0.9472576624920083, 1.768520751856667, 15.91668676671, 0.01023403406619612, 0.01023403406619612),
    @("sCs", "Tnfsf10", "Tnfrsf11b", "sCs", 3, 1, 0.2551576666666667, 0.765473, 0.01080385461255898, 0.01080385461255898, 3, 1, 0.3101076666666667, 0.930323, 0.04238177017576132, 0.04238177017576132, 0.07912634864211111, 0.712137137779, 0.0004578864832018135, 0.0004578864832018135)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $excelRow = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($excelRow, $j + 1).Value = $row[$j]
    }
}

Write-Host "Rows written:" $data.Length
